# The workbook has a single sheet ("Questions") with:
#   A1 = 0                (a leftover numeric placeholder, bold/bordered style)
#   A2 = the "questions = [...]" Python-literal text (shared string)
#
# The target state re-formats that text as pretty-printed JSON (still
# prefixed with "questions = ") and drops the now-unused placeholder
# row, so the JSON text ends up alone in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new, pretty-printed text (single-quoted here-string => fully
# literal, no PowerShell escaping/interpolation needed).
$text = @'
questions = [
    {
        "title": "You\u2019re a business analyst who completed user stories for a web-based software service. Alongside each story, you've included a wireframe.Which purpose should these wireframes serve?",
        "ques_type": 2,
        "options": [
            "To display story content layouts.",
            "To detail design elements.",
            "To follow human-computer interaction guidelines.",
            "To present the user flow."
        ],
        "score": "To present the user flow."
    },
    {
        "title": "You\u2019re a business analyst tracking the pace at which user stories are being marked as done in each sprint.Which of the following factors should influence the number of user stories completed per sprint?",
        "ques_type": 15,
        "options": [
            "The speed and skills of the development team.",
            "The complexity of the user stories.",
            "The number of hours worked by the development team. ",
            "The project completion date.",
            "The capabilities of the project manager."
        ],
        "score": [
            "The speed and skills of the development team.",
            "The complexity of the user stories."
        ]
    },
    {
        "title": "You\u2019re a business analyst who occasionally needs to communicate with the quality assurance (QA) team during the product life cycle.When should you look to do this?",
        "ques_type": 2,
        "options": [
            "Continuously throughout the life cycle",
            "During the user acceptance testing (UAT) phase",
            "During the requirements analysis phase",
            "During the requirements verification and validation phase"
        ],
        "score": "Continuously throughout the life cycle"
    },
    {
        "title": "You are a business analyst and need to evaluate the performance of a recently implemented customer relationship management (CRM) system.Which of the following factors should you prioritize?",
        "ques_type": 2,
        "options": [
            "The system\u2019s behavior during runtime.",
            "The system\u2019s performance in relation to the business goals.",
            "The system\u2019s performance in relation to the validated requirements.",
            "User acceptance tests outcomes."
        ],
        "score": "The system\u2019s performance in relation to the business goals."
    }
]
'@

# Write the reformatted text into A2 first (reuses/updates the existing
# shared-string slot instead of creating a brand-new one).
$ws.Range("A2").Value = $text

# Setting a multi-line value auto-expands the row height; put it back to
# the sheet's normal auto height before we shift rows around.
$ws.Rows("2").AutoFit()

# Remove the old placeholder row (A1 = 0). This shifts row 2 - now
# holding the updated text - up into row 1, matching the target layout.
$ws.Rows("1").Delete()
